$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 73.333336
$ws.Range("I5").Value = 73.333336
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 73.333336
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 41.666664
$ws.Range("N5").ClearContents() | Out-Null
$ws.Range("H11").Value = 118.052635
$ws.Range("I11").Value = 118.052635
$ws.Range("K11").Value = 118.052635
$ws.Range("M11").Value = 21.947365
$ws.Range("H12").Value = 191.04546
$ws.Range("I12").Value = 182
$ws.Range("J12").Value = 221.8
$ws.Range("K12").Value = 182
$ws.Range("L12").Value = 221.8
$ws.Range("M12").Value = -12
$ws.Range("N12").Value = -561.8
$ws.Range("H18").Value = 275.46155
$ws.Range("I18").Value = 192.17392
$ws.Range("J18").Value = 914
$ws.Range("K18").Value = 192.17392
$ws.Range("L18").Value = 914
$ws.Range("M18").Value = 91.82607999999999
$ws.Range("N18").Value = -1482
$ws.Range("H19").Value = 157.1579
$ws.Range("I19").Value = 119.25
$ws.Range("J19").Value = 184.72728
$ws.Range("K19").Value = 119.25
$ws.Range("L19").Value = 184.72728
$ws.Range("M19").Value = 55.75
$ws.Range("N19").Value = -534.7272800000001
$ws.Range("H41").Value = 194.65218
$ws.Range("I41").Value = 200
$ws.Range("J41").Value = 189.75
$ws.Range("K41").Value = 200
$ws.Range("L41").Value = 189.75
$ws.Range("M41").Value = 240
$ws.Range("N41").Value = -1069.75
$ws.Range("H43").Value = 65166.824
$ws.Range("I43").Value = 25237.5
$ws.Range("J43").Value = 77452.766
$ws.Range("K43").Value = 25237.5
$ws.Range("L43").Value = 77452.766
$ws.Range("M43").Value = -25168.5
$ws.Range("N43").Value = -77590.766
$ws.Range("H53").Value = 357.66666
$ws.Range("I53").Value = 396.66666
$ws.Range("J53").Value = 344.66666
$ws.Range("K53").Value = 396.66666
$ws.Range("L53").Value = 344.66666
$ws.Range("M53").Value = 240.33334
$ws.Range("N53").Value = -1618.66666
$ws.Range("H100").Value = 9567.467000000001
$ws.Range("I100").Value = 5277.778
$ws.Range("J100").Value = 16002
$ws.Range("K100").Value = 5277.778
$ws.Range("L100").Value = 16002
$ws.Range("M100").Value = -4736.778
$ws.Range("N100").Value = -17084
$ws.Range("H116").Value = 5083.5244
$ws.Range("I116").Value = 5958.552
$ws.Range("J116").Value = 4290.5312
$ws.Range("K116").Value = 5958.552
$ws.Range("L116").Value = 4290.5312
$ws.Range("M116").Value = -2516.552
$ws.Range("N116").Value = -11174.5312
$ws.Range("H129").Value = 976.4
$ws.Range("J129").Value = 996.25
$ws.Range("L129").Value = 2988.75
$ws.Range("N129").Value = -12988.75
$ws.Range("H137").Value = 4948.84
$ws.Range("I137").Value = 918.5833
$ws.Range("J137").Value = 8669.076999999999
$ws.Range("K137").Value = 2755.7499
$ws.Range("L137").Value = 26007.231
$ws.Range("M137").Value = -205.7498999999998
$ws.Range("N137").Value = -31107.231
$ws.Range("H138").Value = 131723.03
$ws.Range("I138").Value = 2530.7693
$ws.Range("J138").Value = 157561.48
$ws.Range("K138").Value = 7592.3079
$ws.Range("L138").Value = 472684.4400000001
$ws.Range("M138").Value = -2452.3079
$ws.Range("N138").Value = -482964.4400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1039.3462
$ws.Range("I45").Value = 906.4
$ws.Range("J45").Value = 1482.5
$ws.Range("K45").Value = 906.4
$ws.Range("L45").Value = 1482.5
$ws.Range("M45").Value = -529.4
$ws.Range("N45").Value = -2236.5
$ws.Range("H122").Value = 845.88464
$ws.Range("I122").Value = 758.95
$ws.Range("J122").Value = 1135.6666
$ws.Range("K122").Value = 2276.85
$ws.Range("L122").Value = 3406.9998
$ws.Range("M122").Value = 173.1499999999996
$ws.Range("N122").Value = -8306.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2775.2
$ws.Range("I94").Value = 474.57895
$ws.Range("J94").Value = 6749
$ws.Range("K94").Value = 474.57895
$ws.Range("L94").Value = 6749
$ws.Range("M94").Value = -23.57895000000002
$ws.Range("N94").Value = -7651

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 18559.75
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 18559.75
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 18559.75
$ws.Range("M31").ClearContents() | Out-Null
$ws.Range("N31").Value = -19149.75
$ws.Range("H34").Value = 18559.75
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 18559.75
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 18559.75
$ws.Range("M34").ClearContents() | Out-Null
$ws.Range("N34").Value = -18963.75
$ws.Range("H50").Value = 9855
$ws.Range("J50").Value = 10114.546
$ws.Range("L50").Value = 10114.546
$ws.Range("N50").Value = -11364.546
$ws.Range("H51").Value = 9249.786
$ws.Range("J51").Value = 9253.615
$ws.Range("L51").Value = 9253.615
$ws.Range("N51").Value = -10725.615
$ws.Range("H59").Value = 12694.889
$ws.Range("J59").Value = 12694.889
$ws.Range("L59").Value = 12694.889
$ws.Range("N59").Value = -14984.889
$ws.Range("H60").Value = 8343.053
$ws.Range("J60").Value = 8343.053
$ws.Range("L60").Value = 8343.053
$ws.Range("N60").Value = -9365.053
$ws.Range("H61").Value = 9249.786
$ws.Range("J61").Value = 9253.615
$ws.Range("L61").Value = 9253.615
$ws.Range("N61").Value = -9949.615
$ws.Range("H88").Value = 10000
$ws.Range("I88").Value = 10000
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 10000
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -9594
$ws.Range("N88").ClearContents() | Out-Null
$ws.Range("H91").Value = 10000
$ws.Range("I91").Value = 10000
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 10000
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = -8596
$ws.Range("N91").ClearContents() | Out-Null
$ws.Range("H132").Value = 16638.787
$ws.Range("I132").Value = 19723.686
$ws.Range("J132").Value = 2756.75
$ws.Range("K132").Value = 59171.058
$ws.Range("L132").Value = 8270.25
$ws.Range("M132").Value = -56641.058
$ws.Range("N132").Value = -13330.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 27551.75
$ws.Range("J39").Value = 27551.75
$ws.Range("L39").Value = 82655.25
$ws.Range("N39").Value = -83243.25
$ws.Range("H131").Value = 264079.44
$ws.Range("J131").Value = 304019.06
$ws.Range("L131").Value = 912057.1799999999
$ws.Range("N131").Value = -922137.1799999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 18566.56
$ws.Range("I102").Value = 12560.333
$ws.Range("J102").Value = 24110.77
$ws.Range("K102").Value = 12560.333
$ws.Range("L102").Value = 24110.77
$ws.Range("M102").Value = -10938.333
$ws.Range("N102").Value = -27354.77
$ws.Range("H122").Value = 2102.0322
$ws.Range("I122").Value = 2102.5217
$ws.Range("J122").Value = 2100.625
$ws.Range("K122").Value = 6307.5651
$ws.Range("L122").Value = 6301.875
$ws.Range("M122").Value = -3857.5651
$ws.Range("N122").Value = -11201.875
$ws.Range("H132").Value = 18580.967
$ws.Range("I132").Value = 1653.9762
$ws.Range("J132").Value = 58077.277
$ws.Range("K132").Value = 4961.9286
$ws.Range("L132").Value = 174231.831
$ws.Range("M132").Value = -2431.9286
$ws.Range("N132").Value = -179291.831

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2227.5
$ws.Range("I93").Value = 1790.3846
$ws.Range("J93").Value = 2858.889
$ws.Range("K93").Value = 1790.3846
$ws.Range("L93").Value = 2858.889
$ws.Range("M93").Value = -542.3846000000001
$ws.Range("N93").Value = -5354.889
$ws.Range("H136").Value = 418685.1
$ws.Range("I136").Value = 668504.7
$ws.Range("J136").Value = 2319.111
$ws.Range("K136").Value = 2005514.1
$ws.Range("L136").Value = 6957.333
$ws.Range("M136").Value = -2002964.1
$ws.Range("N136").Value = -12057.333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 396.8421
$ws.Range("I113").Value = 339.13333
$ws.Range("J113").Value = 613.25
$ws.Range("K113").Value = 1017.39999
$ws.Range("L113").Value = 1839.75
$ws.Range("M113").Value = 1152.60001
$ws.Range("N113").Value = -6179.75
$ws.Range("H122").Value = 2066.842
$ws.Range("I122").Value = 1176.4286
$ws.Range("J122").Value = 4560
$ws.Range("K122").Value = 3529.2858
$ws.Range("L122").Value = 13680
$ws.Range("M122").Value = -1079.2858
$ws.Range("N122").Value = -18580
